$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Insert two new blank rows after row 113 (new rows 114 & 115), pushing the
#    rest of the "category" section down by two rows (old row 114 -> 116, ...).
# ---------------------------------------------------------------------------
$ws.Rows("114:115").Insert()

# The freshly inserted rows don't pick up the exact blank-row formatting used
# throughout this section (green fill + thin border + date format on column
# C), so copy it over explicitly from row 113, which already has it.
$ws.Range("B113:C113").Copy()
$ws.Range("B114:C114").PasteSpecial(-4122)
$ws.Range("B113:C113").Copy()
$ws.Range("B115:C115").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 2. Populate the new wish-list entry that now lives on row 112 (it was a
#    blank - but already styled - row before).
# ---------------------------------------------------------------------------
$ws.Range("A112").Value = "DONE"
$ws.Range("D112").Value = "robot draw order  during jumps needs to be adjusted"

# Match B112's styling (fill/border, general format) to the rest of the row
# by copying the format from A111 (a cell that already uses that exact style)
$ws.Range("A111").Copy()
$ws.Range("B112").PasteSpecial(-4122)
$ws.Range("B112").Value = "Dave/Tom"
$ws.Range("C112").Value = 39964

# ---------------------------------------------------------------------------
# 3. Row 113 loses its formatted-but-empty B cell - only C113 stays styled.
# ---------------------------------------------------------------------------
$ws.Range("B113").Clear()

# ---------------------------------------------------------------------------
# 4. Make sure the sheet's recorded dimension grows to match the two
#    newly-added rows (touch the last blank separator row so the used range
#    is recalculated, then drop the value back out again).
# ---------------------------------------------------------------------------
$ws.Range("D159").Value = "x"
$ws.Range("D159").ClearContents()

# ---------------------------------------------------------------------------
# 5. Update the view's selection to match (B113).
# ---------------------------------------------------------------------------
$ws.Activate()
$ws.Range("B113").Select()
